# Insert 4 new weekly-report rows (Especial/Primera/Segunda/Tercera) for the
# new price date 2023-11-06 (serial 45236) above the existing Frutilla data
# block, pushing the previously-existing rows down by 4 (dimension grows
# from A1:T149 to A1:T153).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything at/after row 88 down by 4 rows, carrying formatting
# (including the date-number style on column D) down with it.
$ws.Rows("88:91").Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "$/bandeja 3 kilos"
$origen      = "Región de Arica y Parinacota"
$kgUnidad    = 3

$fecha = 45236

$filas = @(
    @{ Row = 88; Calidad = "Especial"; Volumen = 180; PMin = 7000; PMax = 8000; PProm = 7444; PKg = 2481 },
    @{ Row = 89; Calidad = "Primera";  Volumen = 280; PMin = 5000; PMax = 6000; PProm = 5536; PKg = 1845 },
    @{ Row = 90; Calidad = "Segunda";  Volumen = 310; PMin = 3000; PMax = 4000; PProm = 3516; PKg = 1172 },
    @{ Row = 91; Calidad = "Tercera";  Volumen = 160; PMin = 2000; PMax = 3000; PProm = 2625; PKg = 875 }
)

foreach ($f in $filas) {
    $r = $f.Row
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $f.Calidad
    $ws.Cells.Item($r, 13).Value = $f.Volumen
    $ws.Cells.Item($r, 14).Value = $f.PMin
    $ws.Cells.Item($r, 15).Value = $f.PMax
    $ws.Cells.Item($r, 16).Value = $f.PProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $f.PKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
